$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the existing row 108, pushing the
# current rows 109:179 down to 111:181 (dimension becomes A1:T181).
$ws.Rows.Item(109).Insert()
$ws.Rows.Item(109).Insert()

# --- New row 109 ---
$ws.Cells.Item(109, 1).Value = 1
$ws.Cells.Item(109, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(109, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(109, 4).Value = 45086
$ws.Cells.Item(109, 5).Value = 15
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100108
$ws.Cells.Item(109, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(109, 9).Value = 100108003
$ws.Cells.Item(109, 10).Value = "Maracuyá"
$ws.Cells.Item(109, 11).Value = "Sin especificar"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 150
$ws.Cells.Item(109, 14).Value = 34000
$ws.Cells.Item(109, 15).Value = 35000
$ws.Cells.Item(109, 16).Value = 34400
$ws.Cells.Item(109, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(109, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(109, 19).Value = 1720
$ws.Cells.Item(109, 20).Value = 20

# --- New row 110 ---
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(110, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(110, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(110, 4).Value = 45086
$ws.Cells.Item(110, 5).Value = 15
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100108
$ws.Cells.Item(110, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(110, 9).Value = 100108003
$ws.Cells.Item(110, 10).Value = "Maracuyá"
$ws.Cells.Item(110, 11).Value = "Sin especificar"
$ws.Cells.Item(110, 12).Value = "Segunda"
$ws.Cells.Item(110, 13).Value = 160
$ws.Cells.Item(110, 14).Value = 28000
$ws.Cells.Item(110, 15).Value = 30000
$ws.Cells.Item(110, 16).Value = 29250
$ws.Cells.Item(110, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(110, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 19).Value = 1462
$ws.Cells.Item(110, 20).Value = 20
